$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price observation is inserted as row 12 (Región de La
# Araucanía, $/docena de atados (12 kilos)); every existing record from
# the old row 12 onward shifts down by one row.
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C12").Value = 'Los Lagos'
$ws.Range("D12").Value = 44545
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112009
$ws.Range("G12").Value = 'Acelga'
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = '$/docena de atados (12 kilos)'
$ws.Range("O12").Value = 'Región de La Araucanía'
$ws.Range("P12").Value = 833
$ws.Range("Q12").Value = 12
$ws.Range("R12").Value = 'Hortaliza'
